# Rebuild the "Clientes" sheet with the new, smaller column layout.
#
# Obsolete columns removed: tramite, entidad, monto_alcanza, plazo,
# estado_civil, tipo_vivienda, tiempo_pensionado, contrasena_sipre,
# ref1_nombre, ref1_telefono, ref1_parentesco, ref2_nombre, ref2_telefono,
# ref2_parentesco, asesor_venta, fuente_base_nombre, fecha_proximo (incl.
# every trace of the old "segundo_estatus"-era schema).
#
# New columns added: sucursal, asesor, monto_propuesta, monto_final,
# analista. The sample/integration-test row is refreshed with new values
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header row (A1:O1) ------------------------------------------------
# Overwrite header labels in place (columns A-O already carry the bold /
# bordered / centered header style from the source workbook, so we leave
# that formatting untouched rather than clearing + reapplying it).
$headers = @(
    "id",
    "nombre",
    "sucursal",
    "asesor",
    "fecha_ingreso",
    "fecha_dispersion",
    "estatus",
    "monto_propuesta",
    "monto_final",
    "observaciones",
    "score",
    "telefono",
    "correo",
    "analista",
    "fuente"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Drop every column beyond the new O1:O2 boundary (old P1:AC2 range held the
# now-removed fields) so the used range / dimension shrinks back down.
$ws.Range("P1:AC2").Clear()

# ---- New data / test row (row 2) -------------------------------------------
# Plain text values - assign directly so they land as shared strings without
# picking up a text number format / extra style.
$ws.Range("A2").Value = "C1000"
$ws.Range("B2").Value = "Cliente Integración"
$ws.Range("C2").Value = "TOXQUI"
$ws.Range("D2").Value = "Asesor Test"
$ws.Range("G2").Value = "PENDIENTE CLIENTE"
$ws.Range("J2").Value = "Cliente creado por test de integración"
$ws.Range("M2").Value = "integracion@test.com"
$ws.Range("N2").Value = "Test Analyst"
$ws.Range("O2").Value = "Test"

# Numeric-looking values (a date string and plain integers) must stay as
# literal text rather than be reinterpreted as a date serial / number, so
# force text format before assigning, then drop back to the default style.
$textCells = @("E2", "H2", "K2", "L2")
$textValues = @("2025-12-10", "50000", "700", "5551234567")
for ($i = 0; $i -lt $textCells.Length; $i++) {
    $cell = $ws.Range($textCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$i]
    $cell.Style = "Normal"
}

# F2 (fecha_dispersion) and I2 (monto_final) are intentionally left blank.
